# Finished extracting player data, next step is rank data
$wb = $excel.ActiveWorkbook

# Rename the "game_types" sheet to "rank"
$rankSheet = $wb.Worksheets.Item("game_types")
$rankSheet.Name = "rank"

# Make "rank" the active sheet/tab (was "players") - this also clears
# tabSelected on "players" and sets it on "rank"
$rankSheet.Activate()

# Update the selection on the rank sheet to span the header row A1:G1
# (anchored/active at G1, the last header cell).
$rankSheet.Range("A1:G1").Select()
